# Auto-generated edit script applying the Leviathan_Profits.xlsx diff
# Updates currentAveragePrice / Leve price / profit columns (H-N) across several crafting-job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3831.5
$ws.Range("J19").Value = 5248
$ws.Range("L19").Value = 5248
$ws.Range("N19").Value = -5598
$ws.Range("H32").Value = 5584.125
$ws.Range("I32").Value = 5500
$ws.Range("J32").Value = 5634.6
$ws.Range("K32").Value = 5500
$ws.Range("L32").Value = 5634.6
$ws.Range("M32").Value = -5174
$ws.Range("N32").Value = -6286.6
$ws.Range("H53").Value = 397.9
$ws.Range("I53").Value = 231
$ws.Range("K53").Value = 231
$ws.Range("M53").Value = 406
$ws.Range("H62").Value = 95679.55
$ws.Range("I62").Value = 115053.336
$ws.Range("K62").Value = 115053.336
$ws.Range("M62").Value = -114429.336
$ws.Range("H65").Value = 95679.55
$ws.Range("I65").Value = 115053.336
$ws.Range("K65").Value = 575266.6799999999
$ws.Range("M65").Value = -572146.6799999999
$ws.Range("H137").Value = 4362.486
$ws.Range("I137").Value = 1522.9333
$ws.Range("K137").Value = 4568.7999
$ws.Range("M137").Value = -2018.7999
$ws.Range("H138").Value = 2245.8462
$ws.Range("I138").Value = 1369.55
$ws.Range("J138").Value = 3168.2632
$ws.Range("K138").Value = 4108.65
$ws.Range("L138").Value = 9504.7896
$ws.Range("M138").Value = 1031.35
$ws.Range("N138").Value = -19784.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 83054
$ws.Range("I32").Value = 52670.75
$ws.Range("K32").Value = 52670.75
$ws.Range("M32").Value = -52383.75
$ws.Range("H132").Value = 1259.8302
$ws.Range("I132").Value = 939.55817
$ws.Range("J132").Value = 2637
$ws.Range("K132").Value = 2818.67451
$ws.Range("L132").Value = 7911
$ws.Range("M132").Value = -288.6745099999998
$ws.Range("N132").Value = -12971

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 30354
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 30354
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 30354
$ws.Range("M87").Value = ""
$ws.Range("N87").Value = -32850
$ws.Range("H90").Value = 30354
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 30354
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 91062
$ws.Range("M90").Value = ""
$ws.Range("N90").Value = -103542
$ws.Range("H105").Value = 5773.2354
$ws.Range("J105").Value = 5302
$ws.Range("L105").Value = 5302
$ws.Range("N105").Value = -8796

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1676.55
$ws.Range("I31").Value = 1148.5
$ws.Range("J31").Value = 3788.75
$ws.Range("K31").Value = 1148.5
$ws.Range("L31").Value = 3788.75
$ws.Range("M31").Value = -853.5
$ws.Range("N31").Value = -4378.75
$ws.Range("H34").Value = 1676.55
$ws.Range("I34").Value = 1148.5
$ws.Range("J34").Value = 3788.75
$ws.Range("K34").Value = 1148.5
$ws.Range("L34").Value = 3788.75
$ws.Range("M34").Value = -946.5
$ws.Range("N34").Value = -4192.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 444.6
$ws.Range("I113").Value = 561
$ws.Range("J113").Value = 415.5
$ws.Range("K113").Value = 1683
$ws.Range("L113").Value = 1246.5
$ws.Range("M113").Value = 487
$ws.Range("N113").Value = -5586.5
$ws.Range("H122").Value = 2705
$ws.Range("J122").Value = 2705
$ws.Range("L122").Value = 24345
$ws.Range("N122").Value = -29245
$ws.Range("H131").Value = 101367.84
$ws.Range("I131").Value = 333749.66
$ws.Range("K131").Value = 1001248.98
$ws.Range("M131").Value = -996208.98

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 244.2
$ws.Range("I2").Value = 235.83333
$ws.Range("K2").Value = 235.83333
$ws.Range("M2").Value = -122.83333
$ws.Range("H122").Value = 173938.2
$ws.Range("I122").Value = 180007.42
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 540022.26
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -537572.26
$ws.Range("N122").Value = -16900
$ws.Range("H126").Value = 4030
$ws.Range("I126").Value = 4012
$ws.Range("J126").Value = 4057
$ws.Range("K126").Value = 12036
$ws.Range("L126").Value = 12171
$ws.Range("M126").Value = -9566
$ws.Range("N126").Value = -17111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 856.5
$ws.Range("I16").Value = 908.625
$ws.Range("J16").Value = 648
$ws.Range("K16").Value = 908.625
$ws.Range("L16").Value = 648
$ws.Range("M16").Value = -738.625
$ws.Range("N16").Value = -988
$ws.Range("H93").Value = 1554.3334
$ws.Range("I93").Value = 1189.7894
$ws.Range("J93").Value = 2939.6
$ws.Range("K93").Value = 1189.7894
$ws.Range("L93").Value = 2939.6
$ws.Range("M93").Value = 58.21060000000011
$ws.Range("N93").Value = -5435.6
$ws.Range("H94").Value = 30330
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H122").Value = 28575
$ws.Range("J122").Value = 4766.6665
$ws.Range("L122").Value = 14299.9995
$ws.Range("N122").Value = -19199.9995
$ws.Range("H123").Value = 30500
$ws.Range("I123").Value = 20000
$ws.Range("J123").Value = 41000
$ws.Range("K123").Value = 20000
$ws.Range("L123").Value = 41000
$ws.Range("M123").Value = -15100
$ws.Range("N123").Value = -50800
$ws.Range("H136").Value = 3016.875
$ws.Range("I136").Value = 2248.4443
$ws.Range("J136").Value = 3477.9333
$ws.Range("K136").Value = 6745.3329
$ws.Range("L136").Value = 10433.7999
$ws.Range("M136").Value = -4195.3329
$ws.Range("N136").Value = -15533.7999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 45969
$ws.Range("J46").Value = 45969
$ws.Range("L46").Value = 45969
$ws.Range("N46").Value = -46431
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H113").Value = 1122.3334
$ws.Range("I113").Value = 408.42856
$ws.Range("J113").Value = 1576.6364
$ws.Range("K113").Value = 1225.28568
$ws.Range("L113").Value = 4729.9092
$ws.Range("M113").Value = 944.71432
$ws.Range("N113").Value = -9069.9092
$ws.Range("H122").Value = 5388.4707
$ws.Range("I122").Value = 5106.933
$ws.Range("K122").Value = 15320.799
$ws.Range("M122").Value = -12870.799
$ws.Range("H126").Value = 12534.318
$ws.Range("I126").Value = 14962.177
$ws.Range("K126").Value = 44886.531
$ws.Range("M126").Value = -42416.531
$ws.Range("H132").Value = 45626.47
$ws.Range("I132").Value = 45443.668
$ws.Range("K132").Value = 136331.004
$ws.Range("M132").Value = -133801.004
$ws.Range("H134").Value = 45969
$ws.Range("J134").Value = 45969
$ws.Range("L134").Value = 137907
$ws.Range("N134").Value = -142977
